$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column O (rows 3-14) into column P so the new column
# picks up identical styles to its left neighbor.
$ws.Range("O3:O14").Copy()
[void]$ws.Range("P3:P14").PasteSpecial(-4122) # xlPasteFormats

# Fill in the new 2022 column values
$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 96.969944810665083
$ws.Range("P6").Value = 96.173557859042035
$ws.Range("P7").Value = 62.289845326160055
$ws.Range("P8").Value = 100
$ws.Range("P9").Value = 100
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = 100
$ws.Range("P12").Value = 58.090784503861151
$ws.Range("P13").Value = 100
$ws.Range("P14").Value = 100

$excel.CutCopyMode = 0

[void]$ws.Range("Q4").Select()
